$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "other cows"
$ws.Range("D2").Value = ">9000"

$ws.Range("D3").Select()
